# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (1 -> 3) and all of the
# downstream NATMI-derived statistics that depend on those counts, for the
# Efnb2-Ephb4 ligand-receptor pair rows (rows 2-10 of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ E=3; G=20.94432133333333;  H=62.832964;   I=0.7396577289668299;  J=0.7396577289668298;  K=3; M=17.00372766666667; N=51.011183; O=0.6620200065567142; P=0.6620200065567141; Q=356.1315361151569; R=3205.183825036412; S=0.4896682145803451;  T=0.4896682145803449  }
    3  = @{ E=3; G=20.94432133333333;  H=62.832964;   I=0.7396577289668299;  J=0.7396577289668298;  K=3; M=5.393811333333333; N=16.181434; O=0.2100016587103467; P=0.2100016587103467; Q=112.9697177767084; R=1016.727459990376; S=0.1553293499609623;  T=0.1553293499609623  }
    4  = @{ E=3; G=20.94432133333333;  H=62.832964;   I=0.7396577289668299;  J=0.7396577289668298;  K=3; M=3.287074;          N=9.861222;  O=0.1279783347329392; P=0.1279783347329391; Q=68.84553410244534; R=619.609806922008;  S=0.09466016442552254; T=0.09466016442552251 }
    5  = @{ E=3; G=2.327094666666667;  H=6.981284;    I=0.08218235047311259; J=0.08218235047311258; K=3; M=17.00372766666667; N=51.011183; O=0.6620200065567142; P=0.6620200065567141; Q=39.56928396655245; R=356.123555698972;  S=0.05440636019905618; T=0.05440636019905617 }
    6  = @{ E=3; G=2.327094666666667;  H=6.981284;    I=0.08218235047311259; J=0.08218235047311258; K=3; M=5.393811333333333; N=16.181434; O=0.2100016587103467; P=0.2100016587103467; Q=12.55190958680622; R=112.967186281256;  S=0.01725842991606868; T=0.01725842991606868 }
    7  = @{ E=3; G=2.327094666666667;  H=6.981284;    I=0.08218235047311259; J=0.08218235047311258; K=3; M=3.287074;          N=9.861222;  O=0.1279783347329392; P=0.1279783347329391; Q=7.649332374338668; R=68.843991369048;   S=0.01051756035798772; T=0.01051756035798772 }
    8  = @{ E=3; G=5.044818;           H=15.134454;   I=0.1781599205600575;  J=0.1781599205600575;  K=3; M=17.00372766666667; N=51.011183; O=0.6620200065567142; P=0.6620200065567141; Q=85.780711399898;   R=772.026402599082;  S=0.117945431777313;   T=0.117945431777313   }
    9  = @{ E=3; G=5.044818;           H=15.134454;   I=0.1781599205600575;  J=0.1781599205600575;  K=3; M=5.393811333333333; N=16.181434; O=0.2100016587103467; P=0.2100016587103467; Q=27.210796503004;   R=244.897168527036;  S=0.03741387883331567; T=0.03741387883331568 }
    10 = @{ E=3; G=5.044818;           H=15.134454;   I=0.1781599205600575;  J=0.1781599205600575;  K=3; M=3.287074;          N=9.861222;  O=0.1279783347329392; P=0.1279783347329391; Q=16.582690082532;   R=149.244210742788;  S=0.02280060994942889; T=0.02280060994942889 }
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
